$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Label" header column H (copy formatting from the existing header cell G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("H1").Value = "Label"

# Updated D/E fitted values for rows 2-11 (refit of NCDEs to individual patients)
$dValues = @{
    2  = 0.485734567341533
    3  = 0.5984103639214586
    4  = 0.4413629089234173
    5  = 0.5471933407477756
    6  = 0.6840058374851964
    7  = 0.4443539979217413
    8  = 0.5727673579200213
    9  = 0.6412025785058015
    10 = 0.6379959802449519
    11 = 0.6628165440696163
}
$eValues = @{
    2  = 0.485734567341533
    3  = 0.5984103639214586
    4  = 0.4413629089234173
    5  = 0.5471933407477756
    6  = 0.6840058374851964
    7  = 0.5556460020782588
    8  = 0.4272326420799787
    9  = 0.3587974214941985
    10 = 0.3620040197550481
    11 = 0.3371834559303837
}

foreach ($r in $dValues.Keys) {
    $ws.Cells.Item($r, 4).Value = $dValues[$r]
    $ws.Cells.Item($r, 5).Value = $eValues[$r]
}

# New Label (H) column values for rows 2-21: 0 for Control rows, 1 for MDD rows
$hValues = @{
    2 = 0; 3 = 0; 4 = 0; 5 = 0; 6 = 0; 7 = 1; 8 = 1; 9 = 1; 10 = 1; 11 = 1
    12 = 0; 13 = 0; 14 = 0; 15 = 0; 16 = 0; 17 = 1; 18 = 1; 19 = 1; 20 = 1; 21 = 1
}

foreach ($r in $hValues.Keys) {
    $ws.Cells.Item($r, 8).Value = $hValues[$r]
}
